$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("trading_journal")

# Shift existing header cells right (I1->J1->K1) before writing the new column
$ws.Range("J1").Copy($ws.Range("K1"))
$ws.Range("I1").Copy($ws.Range("J1"))

$ws.Range("I1").Value = "type"
$ws.Range("J1").Value = "end"
$ws.Range("K1").Value = "outcome"
$ws.Range("I2").Value = "short"
